$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff inserts one new price-list row at row 139 ("Segunda" quality,
# 14-unit box, date 44491) and pushes the existing rows 139-160 down to
# 140-161 (the old row 160 ends up at 161 unchanged).
$ws.Rows("139:139").Insert()

$ws.Range("A139").Value = 5
$ws.Range("B139").Value = "Macroferia Regional de Talca"
$ws.Range("C139").Value = "Maule"
$ws.Range("D139").Value = 44491
$ws.Range("E139").Value = 7
$ws.Range("F139").Value = "Fruta"
$ws.Range("G139").Value = 100108
$ws.Range("H139").Value = "Tropicales y subtropicales"
$ws.Range("I139").Value = 100108005
$ws.Range("J139").Value = "Piña"
$ws.Range("K139").Value = "Caramelo"
$ws.Range("L139").Value = "Segunda"
$ws.Range("M139").Value = 50
$ws.Range("N139").Value = 20000
$ws.Range("O139").Value = 20000
$ws.Range("P139").Value = 20000
$ws.Range("Q139").Value = "$/caja 14 unidades"
$ws.Range("R139").Value = "Ecuador"
$ws.Range("S139").Value = 1429
$ws.Range("T139").Value = 14
